$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts existing rows 8..45 down to 9..46)
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(8, 3).Value = "Maule"
$ws.Cells.Item(8, 4).Value = 44707
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 7
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100104
$ws.Cells.Item(8, 8).Value = "Frutos de pepita"
$ws.Cells.Item(8, 9).Value = 100104003
$ws.Cells.Item(8, 10).Value = "Membrillo"
$ws.Cells.Item(8, 11).Value = "Champion"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 120
$ws.Cells.Item(8, 14).Value = 10000
$ws.Cells.Item(8, 15).Value = 10000
$ws.Cells.Item(8, 16).Value = 10000
$ws.Cells.Item(8, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 556
$ws.Cells.Item(8, 20).Value = 18
